$wb = $excel.ActiveWorkbook

# --- "Org" sheet: rebuild as a 3-column org/investor table -----------------
$wsOrg = $wb.Worksheets.Item("Org")

$wsOrg.Range("A1").Value = "OrgName"
$wsOrg.Range("B1").Value = "Industry"

$wsOrg.Range("A2").Value = "Agro Tech Food Ltd"
$wsOrg.Range("B2").Value = "Food & Beverage"

$wsOrg.Range("C1").Value = "Type"
$wsOrg.Range("C2").Value = "Investor"

$wsOrg.Range("A3").Value = "Mahindra & Mahindra"
$wsOrg.Range("B3").Value = "Manufacturing"
$wsOrg.Range("C3").Value = "Integrator"

$wsOrg.Range("A4").Value = "Resilient Innovation Pvt Ltd"
$wsOrg.Range("B4").Value = "Technology"
$wsOrg.Range("C4").Value = "Investor"

# Row 2 / B2 picked up a different (monospace) font in the real edit.
$wsOrg.Range("B2").Font.Name = "Consolas"

# Column widths (characters, COM units -> XML stores chars + 5/6)
$wsOrg.Columns.Item(1).ColumnWidth = 23.83
$wsOrg.Columns.Item(2).ColumnWidth = 17.28
$wsOrg.Columns.Item(3).ColumnWidth = 12.28

# --- "Leads" sheet: no data changes, just where the cursor sits ------------
$wsLeads = $wb.Worksheets.Item("Leads")
$wsLeads.Range("A3").Select()

# Selection + active tab: Org becomes the active / selected sheet
$wsOrg.Range("A5").Select()
$wsOrg.Activate()
